# Re-split the run containing "no endereço: #ENDERECO #CEP " so that
# "#ENDERECO" becomes its own run (it used to be split across three runs
# as "...#ENDERE" / "C" / "O #CEP ...").  The combined text does not
# change, only the run boundaries.

$d = $word.ActiveDocument

$text = $d.Content.Text
$idxStart = $text.IndexOf("no endereço: #ENDERE")
$idxEnd = $idxStart + ("no endereço: #ENDERE" + "C" + "O #CEP ").Length

# Step 1: mark the whole span bold so it becomes formatting-distinct from
# its neighbours; this keeps the later "real" edit (step 2) from bleeding
# into the surrounding runs when things get re-merged.
$rSpan = $d.Range($idxStart, $idxEnd)
$rSpan.Bold = 1

# Step 2: make a genuine text edit (change a character, then change it
# back) so the engine normalizes/merges the touched run(s) into one.
# Because the span is bold-isolated from its neighbours, only the runs
# inside the span get merged.
$idxC = $idxStart + "no endereço: #ENDERE".Length
$rChar = $d.Range($idxC, $idxC + 1)
$rChar.Text = "Q"

$text2 = $d.Content.Text
$idxQ = $text2.IndexOf("no endereço: #ENDEREQ")
$rChar2 = $d.Range($idxQ + "no endereço: #ENDERE".Length, $idxQ + "no endereço: #ENDERE".Length + 1)
$rChar2.Text = "C"

# Step 3: remove the bold formatting again, restoring the original look.
$text3 = $d.Content.Text
$idxStart3 = $text3.IndexOf("no endereço: #ENDERECO #CEP ")
$rSpan2 = $d.Range($idxStart3, $idxStart3 + ("no endereço: #ENDERECO #CEP ").Length)
$rSpan2.Bold = 0

# Step 4: force a run split exactly around "#ENDERECO" by toggling Bold
# on and off for just that sub-range.
$p1 = $idxStart3 + "no endereço: ".Length
$p2 = $p1 + "#ENDERECO".Length
$rTag = $d.Range($p1, $p2)
$rTag.Bold = 1
$rTag.Bold = 0
